# "new sheet for analysis"
#
# Adds a new "analysis" worksheet at the end of the workbook (after
# "pipelines"), populates it with a single data row describing the
# "lichens" repo, and leaves it as the active/selected sheet & tab,
# matching the author's commit.

$wb = $excel.ActiveWorkbook

# Move the previously-selected sheet's cursor off of its old spot onto
# the cell it ends up parked at in the edited workbook.
$wsTutorials = $wb.Worksheets.Item("tutorials")
$wsTutorials.Range("A6").Select()

# Create the new sheet and push it to the very end of the tab strip
# (Worksheets.Add() inserts before the active sheet by default).
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "analysis"
$wsNew.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch by name so we operate on the sheet in its final position.
$ws = $wb.Worksheets.Item("analysis")

# Header row - reuses the same column headers as the other sheets
# (name/description/link/author).
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "link"
$ws.Range("D1").Value = "author"

# Data row for the lichens analysis repo.
$ws.Range("A2").Value = "lichens"
$ws.Range("B2").Value = "Analysis of lichens transcription experiment data from the squamulose R package."
$ws.Range("C2").Value = "https://github.com/TeamMacLean/CJ_NT_1563_23_23022022_lichens/blob/main/README.md"
$ws.Range("D2").Value = "clara"

# Make "analysis" the active sheet/tab and park the selection where the
# author left it.
$ws.Activate()
$ws.Range("D3").Select()

$win = $excel.ActiveWindow
$win.Top = -6320
